$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (D) column cells stay text even though the values look numeric
# (matches the source data, which stores every price as an inline string).

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '65.322.39'
$ws.Range("E2").Value = '  -0.79%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.326.49'
$ws.Range("E3").Value = '  -4.27%  '

# Row 4
$ws.Range("E4").Value = '  +0.01%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '571.84'
$ws.Range("E5").Value = '  -1.70%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.74'
$ws.Range("E6").Value = '  +3.09%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.613'
$ws.Range("E7").Value = '  +2.76%  '

# Row 8
$ws.Range("E8").Value = '  +0.02%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.324.29'
$ws.Range("E9").Value = '  -4.25%  '

# Row 10
$ws.Range("E10").Value = '  -1.86%  '

# Row 11
$ws.Range("E11").Value = '  -0.22%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.405'
$ws.Range("E12").Value = '  -1.06%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.901.59'
$ws.Range("E13").Value = '  -4.26%  '

# Row 14
$ws.Range("E14").Value = '  +0.03%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.32'
$ws.Range("E15").Value = '  -4.80%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '65.356.52'
$ws.Range("E16").Value = '  -0.84%  '

# Row 17
$ws.Range("E17").Value = '  -1.76%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.329.12'
$ws.Range("E18").Value = '  -4.18%  '

# Row 19
$ws.Range("E19").Value = '  -2.99%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.29'
$ws.Range("E20").Value = '  -3.93%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '360.93'
$ws.Range("E21").Value = '  -1.24%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.39'
$ws.Range("E22").Value = '  -4.15%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.10%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.97'
$ws.Range("E24").Value = '  -2.07%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  -3.27%  '

# Row 26
$ws.Range("E26").Value = '  -3.79%  '

# Row 27
$ws.Range("E27").Value = '  -0.56%  '

# Row 28
$ws.Range("E28").Value = '  -0.97%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.17%  '

# Row 30
$ws.Range("E30").Value = '  -1.35%  '

# Row 31
$ws.Range("E31").Value = '  -0.03%  '

# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '22.84'
$ws.Range("E32").Value = '  -4.71%  '

# Row 33
$ws.Range("B33").Value = 'NEARProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.55'
$ws.Range("E33").Value = '  -3.28%  '

# Row 34
$ws.Range("E34").Value = '  -4.51%  '

# Row 35
$ws.Range("E35").Value = '  -6.52%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '159.82'
$ws.Range("E36").Value = '  -0.39%  '

# Row 37
$ws.Range("E37").Value = '  -3.39%  '

# Row 38
$ws.Range("E38").Value = '  -5.19%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.37'
$ws.Range("E39").Value = '  -6.56%  '

# Row 40
$ws.Range("E40").Value = '  -0.64%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.694.37'
$ws.Range("E41").Value = '  -4.38%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.47'
$ws.Range("E42").Value = '  -2.89%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '6.18'
$ws.Range("E43").Value = '  -4.21%  '

# Row 44
$ws.Range("E44").Value = '  -4.09%  '

# Row 45
$ws.Range("B45").Value = 'OKB'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '39.65'
$ws.Range("E45").Value = '  -1.06%  '

# Row 46
$ws.Range("B46").Value = 'Hedera'
$ws.Range("C46").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0662'
$ws.Range("E46").Value = '  -2.50%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '331.87'
$ws.Range("E47").Value = '  +2.55%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '23.75'

# Row 49
$ws.Range("E49").Value = '  -3.73%  '

# Row 50
$ws.Range("E50").Value = '  +1.69%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.999'
$ws.Range("E51").Value = '  +0.05%  '
